# Append September 14-30, 2020 scheduled-vs-actual flight rows (162-178)
# to the "Ark1" sheet, continuing the existing table pattern:
#   col A = date label (text), col B = scheduled, col C = actual,
#   col D = C/B (percentage formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (number formats, fonts, borders) of the last
# existing data row down across the new rows before filling values in,
# so the new cells pick up the same styles (text col A, integer cols
# B/C, percent formula col D) as the rest of the table.
$ws.Range("A161:D161").Copy()
$ws.Range("A162:D178").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$dates = @(
    "2020-09-14","2020-09-15","2020-09-16","2020-09-17","2020-09-18",
    "2020-09-19","2020-09-20","2020-09-21","2020-09-22","2020-09-23",
    "2020-09-24","2020-09-25","2020-09-26","2020-09-27","2020-09-28",
    "2020-09-29","2020-09-30"
)
$scheduled = @(50,44,53,52,52,43,44,50,50,52,57,58,47,44,51,56,58)
$actual    = @(50,42,51,50,50,42,44,48,49,51,53,57,44,41,30,50,53)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 162 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $scheduled[$i]
    $ws.Cells.Item($r, 3).Value = $actual[$i]
}

# Fill column D with the same relative formula as the rest of the
# table; assigning identical formula text across the whole range in
# one shot lets Excel store it as a shared formula, matching the
# existing D133:D161 shared-formula pattern.
$ws.Range("D162:D178").Formula = "=C162/B162"

# Reflect the author's final scroll position / active cell.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 142
$aw.ScrollColumn = 1
$ws.Range("I173").Select() | Out-Null

"Appended rows 162:178 (2020-09-14 through 2020-09-30)."
